$d = $word.ActiveDocument

$d.Content.Find.Execute("232÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "516÷8=", 2) | Out-Null
$d.Content.Find.Execute("435÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "808÷5=", 2) | Out-Null
$d.Content.Find.Execute("920÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "994÷3=", 2) | Out-Null
$d.Content.Find.Execute("206÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "865÷8=", 2) | Out-Null
$d.Content.Find.Execute("987÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "273÷6=", 2) | Out-Null
$d.Content.Find.Execute("149÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "368÷8=", 2) | Out-Null
$d.Content.Find.Execute("637÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "400÷9=", 2) | Out-Null
$d.Content.Find.Execute("125÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "196÷5=", 2) | Out-Null
$d.Content.Find.Execute("857÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "755÷4=", 2) | Out-Null
$d.Content.Find.Execute("461÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "531÷3=", 2) | Out-Null
$d.Content.Find.Execute("665÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "115÷4=", 2) | Out-Null
$d.Content.Find.Execute("251÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "897÷2=", 2) | Out-Null
$d.Content.Find.Execute("436÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "366÷7=", 2) | Out-Null
$d.Content.Find.Execute("232÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "696÷9=", 2) | Out-Null
$d.Content.Find.Execute("398÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "284÷4=", 2) | Out-Null
$d.Content.Find.Execute("323÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "897÷9=", 2) | Out-Null
$d.Content.Find.Execute("887÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "105÷5=", 2) | Out-Null
$d.Content.Find.Execute("655÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "677÷2=", 2) | Out-Null
$d.Content.Find.Execute("122÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "860÷2=", 2) | Out-Null
$d.Content.Find.Execute("956÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "368÷3=", 2) | Out-Null
$d.Content.Find.Execute("367÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "305÷6=", 2) | Out-Null
$d.Content.Find.Execute("456÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "445÷3=", 2) | Out-Null
$d.Content.Find.Execute("380÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "837÷5=", 2) | Out-Null
$d.Content.Find.Execute("963÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "742÷6=", 2) | Out-Null
$d.Content.Find.Execute("671÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "194÷6=", 2) | Out-Null
